$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 7 / column C: replace the "Valuelist" instruction text for the
# codelist-based TDF_Type with an expanded explanation that documents the
# ct-id / term-weight-list|dset codelist() syntax. The appended
# "term-weight-list|dset" fragment is written in italics to call out the
# placeholder tokens, matching the rest of the workbook's convention.
# ---------------------------------------------------------------------------
$newText = "(contains the ""|""-separated list of values)" + "`n" + "codelist(ct-id, term-weight-list|dset)"
$ws.Range("C7").Value = $newText
$ws.Range("C7").Characters(61, 22).Font.Italic = $true
$ws.Range("C7").Characters(82, 1).Font.Italic = $false

# ---------------------------------------------------------------------------
# Row 7 / column E: add a worked "example" note, highlighted with a yellow
# fill so it stands out next to the new instruction text.
# ---------------------------------------------------------------------------
$ws.Range("E7").Value = "example"
$ws.Range("E7").Interior.Color = 65535
$ws.Range("E7").Font.Italic = $true
$ws.Range("E7").Font.Italic = $false

# ---------------------------------------------------------------------------
# Refresh the active selection/scroll position to rest on the cell that was
# just edited (mirrors the author re-selecting C7 after the call).
# ---------------------------------------------------------------------------
[void]$ws.Range("C7").Select()
